$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table cells to match the latest scrape.
# D-column values that are plain decimal numbers (e.g. "243.52") would be
# auto-coerced to numeric by Excel on assignment, but the source data stores
# them as text (prices like "30.860.42" use "." as a thousands separator, so
# the whole column is kept as text for consistent formatting). Force those
# cells to Text before assigning, then clear the format again so no extra
# cell style lingers on the cell.

$ws.Range("D2").Value = "30.860.42"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.940.90"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4911"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2940"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06893"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.27"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "105.62"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "1.940.63"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07766"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.369"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7051"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "30.901.00"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007729"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.219.80"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.640"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.11"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.555"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.809"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.159"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1040"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.391"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.559"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.578"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.385"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04877"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7573"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.153"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9998"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.737"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02001"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.73"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.488"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.094"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9136"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4437"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.679"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "990.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.69%  "
